# functioneaza comutarea AddGoal / EditGoal. lucrez la EditGoal
#
# The sheet holds a small "Goal/Task" list. Each Goal block reserves a
# blank spacer row below it (for a future task row) and an "Estimated
# Date" cell in column C formatted as dd/mm/yyyy. This edit:
#  - removes the last (experimental/duplicate) "Goal" block that used
#    placeholder text "a" (rows 19-22), so that string falls out of the
#    shared-string table entirely,
#  - stamps every Estimated-Date cell (including the blank spacer rows)
#    with the same "now" timestamp using the dd/mm/yyyy date format,
#  - normalizes the Goal row's Label value (G16/G18) to a plain 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing duplicate Goal entries (rows 20 and 22, plus the
# blank spacer rows 19/21 around them) - only two real Goal blocks with
# the "00:00/0%" template remain afterwards. Select the whole 19:28 band
# first so that, like in Excel, the post-delete selection collapses onto
# the vacated rows.
$ws.Rows("19:28").Select()
$ws.Rows("19:22").Delete()

$now = 44236.456396574074

# Make sure the custom date format used across the sheet reads dd/mm/yyyy
# (lowercase month specifier).
$dateFormat = "dd/mm/yyyy"

# Estimated-Date column for every row that belongs to a Goal/Task block,
# including the blank spacer rows that separate the blocks.
$dateRows = 3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18
foreach ($r in $dateRows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.NumberFormat = $dateFormat
}

# Rows that previously had no Estimated-Date / were blank spacer rows
# stay empty values, only the format is stamped; every other row gets
# "now" written into it.
$valueRows = 3,4,5,6,7,9,10,11,13,14,16,18
foreach ($r in $valueRows) {
    $ws.Cells.Item($r, 3).Value = $now
}

# The "Label" value on the two remaining Goal rows is a plain number.
$ws.Range("G16").Value = 0
$ws.Range("G18").Value = 0
